$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.566.65"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "2.736.63"
$ws.Range("E3").Value = "  +4.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("D15").Value = "3.170.24"
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").Value = "2.746.68"
$ws.Range("E16").Value = "  +5.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.876"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "51.519.36"
$ws.Range("E18").Value = "  +6.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0816"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("E41").Value = "  +10.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.43%  "
$ws.Range("D46").Value = "2.087.08"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
